$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Constants")

# B18 used to hold the display text "Adobe Sign" for the Adobe Sign hyperlink;
# it now shows the literal URL text instead.
$ws.Range("B18").Value = "https://omesofm.na1.documents.adobe.com/account/homeJS"
$ws.Range("B18").Style = "Hyperlink"

# New rows of configuration added below the existing Constants.
$ws.Range("A25").Value = "FillingFieldURL"
$ws.Range("B25").Value = "https://omesofm.na1.documents.adobe.com/account"

$ws.Range("A26").Value = "MinimumDelay"
$ws.Range("B26").Value = 2

$ws.Range("A27").Value = "Sharepoint URL"
$ws.Range("B27").Value = "https://officemgmtentserv.sharepoint.com/sites/ACOE_Automations_DEV"

$ws.Range("A28").Value = "Root Folder"
$ws.Range("B28").Value = "P003_090_TimesheetApprovals"

$ws.Range("A29").Value = "LocalDownloadPath"
$ws.Range("B29").Value = "Data\Input"

# Move the Adobe Sign hyperlink from B18 onto the new Sharepoint URL row (B27).
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B27"), "https://omesofm.na1.documents.adobe.com/account/homeJS")
$ws.Range("B27").Style = "Hyperlink"

# Restore view state: scrolled down a bit with C29 selected.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C29").Select()
